$d = $word.ActiveDocument

$pairs = @(
    @("989÷6=", "561÷6="),
    @("786÷7=", "138÷6="),
    @("228÷3=", "871÷9="),
    @("291÷6=", "422÷6="),
    @("744÷4=", "738÷5="),
    @("922÷8=", "879÷6="),
    @("649÷8=", "182÷4="),
    @("710÷6=", "293÷9="),
    @("178÷9=", "754÷6="),
    @("977÷2=", "836÷2="),
    @("784÷8=", "945÷6="),
    @("784÷4=", "434÷8="),
    @("210÷7=", "524÷9="),
    @("531÷2=", "900÷8="),
    @("803÷6=", "762÷6="),
    @("491÷6=", "117÷3="),
    @("182÷8=", "567÷4="),
    @("946÷2=", "366÷9="),
    @("260÷2=", "263÷6="),
    @("679÷8=", "209÷9="),
    @("174÷2=", "872÷7="),
    @("769÷9=", "698÷2="),
    @("939÷8=", "246÷5="),
    @("113÷6=", "340÷8="),
    @("121÷8=", "633÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
